$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 303.19565
$ws.Range("I8").Value = 277.4
$ws.Range("K8").Value = 832.1999999999999
$ws.Range("M8").Value = -693.1999999999999
$ws.Range("H51").Value = 67422
$ws.Range("J51").Value = 109465.664
$ws.Range("L51").Value = 109465.664
$ws.Range("N51").Value = -110433.664
$ws.Range("H86").Value = 73101920
$ws.Range("I86").Value = 95240024
$ws.Range("J86").Value = 11115230
$ws.Range("K86").Value = 95240024
$ws.Range("L86").Value = 11115230
$ws.Range("M86").Value = -95238901
$ws.Range("N86").Value = -11117476
$ws.Range("H89").Value = 73101920
$ws.Range("I89").Value = 95240024
$ws.Range("J89").Value = 11115230
$ws.Range("K89").Value = 476200120
$ws.Range("L89").Value = 55576150
$ws.Range("M89").Value = -476194504
$ws.Range("N89").Value = -55587382
$ws.Range("H112").Value = 6416.7646
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 6416.7646
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 19250.2938
$ws.Range("M112").Value = ""
$ws.Range("N112").Value = -21466.2938
$ws.Range("H137").Value = 2954.75
$ws.Range("I137").Value = 3811.5715
$ws.Range("K137").Value = 11434.7145
$ws.Range("M137").Value = -8884.7145

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2303.4482
$ws.Range("I2").Value = 1489.4667
$ws.Range("J2").Value = 3175.5715
$ws.Range("K2").Value = 1489.4667
$ws.Range("L2").Value = 3175.5715
$ws.Range("M2").Value = -1376.4667
$ws.Range("N2").Value = -3401.5715
$ws.Range("H4").Value = 197.25
$ws.Range("J4").Value = 30
$ws.Range("L4").Value = 30
$ws.Range("N4").Value = -262
$ws.Range("H5").Value = 211.66667
$ws.Range("I5").Value = 194
$ws.Range("J5").Value = 300
$ws.Range("K5").Value = 194
$ws.Range("L5").Value = 300
$ws.Range("M5").Value = -82
$ws.Range("N5").Value = -524
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").Value = ""
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = ""
$ws.Range("N10").Value = ""
$ws.Range("H116").Value = 2303.4482
$ws.Range("I116").Value = 1489.4667
$ws.Range("J116").Value = 3175.5715
$ws.Range("K116").Value = 1489.4667
$ws.Range("L116").Value = 3175.5715
$ws.Range("M116").Value = 804.5333000000001
$ws.Range("N116").Value = -7763.5715
$ws.Range("H126").Value = 5245
$ws.Range("I126").Value = 5245
$ws.Range("K126").Value = 15735
$ws.Range("M126").Value = -13265
$ws.Range("H132").Value = 6131.7334
$ws.Range("I132").Value = 4839.0625
$ws.Range("J132").Value = 9313.691999999999
$ws.Range("K132").Value = 14517.1875
$ws.Range("L132").Value = 27941.076
$ws.Range("M132").Value = -11987.1875
$ws.Range("N132").Value = -33001.076

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2303.4482
$ws.Range("I3").Value = 1489.4667
$ws.Range("J3").Value = 3175.5715
$ws.Range("K3").Value = 1489.4667
$ws.Range("L3").Value = 3175.5715
$ws.Range("M3").Value = -1375.4667
$ws.Range("N3").Value = -3403.5715
$ws.Range("H4").Value = 211.66667
$ws.Range("I4").Value = 194
$ws.Range("J4").Value = 300
$ws.Range("K4").Value = 194
$ws.Range("L4").Value = 300
$ws.Range("M4").Value = -79
$ws.Range("N4").Value = -530
$ws.Range("H20").Value = 13890871
$ws.Range("J20").Value = 1260.25
$ws.Range("L20").Value = 1260.25
$ws.Range("N20").Value = -1754.25
$ws.Range("H94").Value = 2616.8696
$ws.Range("I94").Value = 1119.1333
$ws.Range("J94").Value = 5425.125
$ws.Range("K94").Value = 1119.1333
$ws.Range("L94").Value = 5425.125
$ws.Range("M94").Value = -668.1333
$ws.Range("N94").Value = -6327.125
$ws.Range("H105").Value = 3691.8572
$ws.Range("I105").Value = 2316.3333
$ws.Range("J105").Value = 4723.5
$ws.Range("K105").Value = 2316.3333
$ws.Range("L105").Value = 4723.5
$ws.Range("M105").Value = -569.3332999999998
$ws.Range("N105").Value = -8217.5
$ws.Range("H113").Value = 5035.5
$ws.Range("I113").Value = 5035.5
$ws.Range("K113").Value = 5035.5
$ws.Range("M113").Value = -2865.5
$ws.Range("H128").Value = 3833
$ws.Range("I128").Value = 3833
$ws.Range("K128").Value = 11499
$ws.Range("M128").Value = -9009
$ws.Range("H132").Value = 103446.336
$ws.Range("J132").Value = 103446.336
$ws.Range("L132").Value = 103446.336
$ws.Range("N132").Value = -113566.336
$ws.Range("H134").Value = 6091.171
$ws.Range("I134").Value = 2693.2104
$ws.Range("J134").Value = 9025.772000000001
$ws.Range("K134").Value = 8079.6312
$ws.Range("L134").Value = 27077.316
$ws.Range("M134").Value = -5544.6312
$ws.Range("N134").Value = -32147.316
$ws.Range("H140").Value = 69853
$ws.Range("J140").Value = 69853
$ws.Range("L140").Value = 69853
$ws.Range("N140").Value = -80213

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H76").Value = 4953.8335
$ws.Range("I76").Value = 4953.8335
$ws.Range("K76").Value = 4953.8335
$ws.Range("M76").Value = -4638.8335
$ws.Range("H79").Value = 4953.8335
$ws.Range("I79").Value = 4953.8335
$ws.Range("K79").Value = 4953.8335
$ws.Range("M79").Value = -3861.8335
$ws.Range("H105").Value = 6498035.5
$ws.Range("I105").Value = 14288215
$ws.Range("J105").Value = 6218.5
$ws.Range("K105").Value = 14288215
$ws.Range("L105").Value = 6218.5
$ws.Range("M105").Value = -14286468
$ws.Range("N105").Value = -9712.5
$ws.Range("H134").Value = 8083.8
$ws.Range("I134").Value = 4353.4614
$ws.Range("J134").Value = 10936.412
$ws.Range("K134").Value = 13060.3842
$ws.Range("L134").Value = 32809.236
$ws.Range("M134").Value = -10525.3842
$ws.Range("N134").Value = -37879.236

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 213.81818
$ws.Range("I23").Value = 207
$ws.Range("J23").Value = 222
$ws.Range("K23").Value = 621
$ws.Range("L23").Value = 666
$ws.Range("M23").Value = -386
$ws.Range("N23").Value = -1136
$ws.Range("H29").Value = 525
$ws.Range("I29").Value = 42.5
$ws.Range("K29").Value = 127.5
$ws.Range("M29").Value = 149.5
$ws.Range("H131").Value = 1773.4706
$ws.Range("J131").Value = 2590.2222
$ws.Range("L131").Value = 7770.6666
$ws.Range("N131").Value = -17850.6666
$ws.Range("H132").Value = 8779.861999999999
$ws.Range("I132").Value = 3477.9375
$ws.Range("K132").Value = 31301.4375
$ws.Range("M132").Value = -28771.4375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 6936
$ws.Range("I102").Value = 5703.7856
$ws.Range("K102").Value = 5703.7856
$ws.Range("M102").Value = -4081.7856
$ws.Range("H122").Value = 2686150.5
$ws.Range("I122").Value = 4263593.5
$ws.Range("J122").Value = 4497.8
$ws.Range("K122").Value = 12790780.5
$ws.Range("L122").Value = 13493.4
$ws.Range("M122").Value = -12788330.5
$ws.Range("N122").Value = -18393.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3103.5
$ws.Range("I22").Value = 1500
$ws.Range("J22").Value = 3504.375
$ws.Range("K22").Value = 1500
$ws.Range("L22").Value = 3504.375
$ws.Range("M22").Value = -1205
$ws.Range("N22").Value = -4094.375
$ws.Range("H27").Value = 3103.5
$ws.Range("I27").Value = 1500
$ws.Range("J27").Value = 3504.375
$ws.Range("K27").Value = 1500
$ws.Range("L27").Value = 3504.375
$ws.Range("M27").Value = -1393
$ws.Range("N27").Value = -3718.375
$ws.Range("H40").Value = 6376.4116
$ws.Range("I40").Value = 3665.8333
$ws.Range("K40").Value = 3665.8333
$ws.Range("M40").Value = -3529.8333
$ws.Range("H46").Value = 3384.2
$ws.Range("J46").Value = 4007
$ws.Range("L46").Value = 4007
$ws.Range("N46").Value = -4383
$ws.Range("H93").Value = 5553.125
$ws.Range("I93").Value = 5203.846
$ws.Range("K93").Value = 5203.846
$ws.Range("M93").Value = -3955.846
$ws.Range("H122").Value = 4637.5835
$ws.Range("I122").Value = 4026.28
$ws.Range("K122").Value = 12078.84
$ws.Range("M122").Value = -9628.84
$ws.Range("H139").Value = 89244.5
$ws.Range("J139").Value = 89244.5
$ws.Range("L139").Value = 89244.5
$ws.Range("N139").Value = -99524.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 5907.727
$ws.Range("I126").Value = 4996.5
$ws.Range("J126").Value = 6110.222
$ws.Range("K126").Value = 14989.5
$ws.Range("L126").Value = 18330.666
$ws.Range("M126").Value = -12519.5
$ws.Range("N126").Value = -23270.666
$ws.Range("H139").Value = 86745.5
$ws.Range("J139").Value = 86094.60000000001
$ws.Range("L139").Value = 86094.60000000001
$ws.Range("N139").Value = -96374.60000000001
